# Initial commit v2p1, R2019a
#
# Re-points H13/H22 "1/0.2" sensitivity formulas to "1/0.025" on the three
# populated sheets, moves the active tab from "None" to "Bus_Makhulu", and
# updates the saved cursor/selection position on each sheet's frozen
# bottom-right pane to match the last-edited cell before save.

$wb = $excel.ActiveWorkbook

$sedanHambaLG = $wb.Worksheets.Item("Sedan_HambaLG")
$sedanHamba   = $wb.Worksheets.Item("Sedan_Hamba")
$busMakhulu   = $wb.Worksheets.Item("Bus_Makhulu")
$none         = $wb.Worksheets.Item("None")

# --- Sedan_HambaLG: update the 1/0.2 -> 1/0.025 formulas -------------------
$sedanHambaLG.Range("H13").Formula = "=1/0.025"
$sedanHambaLG.Range("H22").Formula = "=1/0.025"

# --- Sedan_Hamba: update the 1/0.2 -> 1/0.025 formulas ---------------------
$sedanHamba.Range("H13").Formula = "=1/0.025"
$sedanHamba.Range("H22").Formula = "=1/0.025"

# --- Bus_Makhulu: update the 1/0.2 -> 1/0.025 formulas ---------------------
$busMakhulu.Range("H13").Formula = "=1/0.025"
$busMakhulu.Range("H22").Formula = "=1/0.025"

# --- Update each sheet's remembered selection (bottom-right frozen pane) ---
$sedanHambaLG.Activate()
$sedanHambaLG.Range("H22").Select()

$sedanHamba.Activate()
$sedanHamba.Range("H22").Select()

$busMakhulu.Activate()
$busMakhulu.Range("J22").Select()

# Bus_Makhulu ends up the active/selected tab (was "None").
$busMakhulu.Activate()
